{"js": "// Apply the per-label updates (NO / NAMA / SEPATU size / KELAS) to the\n// three shoe-label cells in the single-row table.\n//\n// Cell 1 (NO : 1): 1 -> T1, ABU NAWAS -> GEDE SUDIANTARA, 42 -> 43,\n//                  DP 4 NAUTIKA / 31 -> DP 5 NAUTIKA / 38\n// Cell 2 (NO : 2): 2 -> T2, ACHMAD UBAIDILLAH -> AGUS SETIAWAN, 43 -> 42,\n//                  DP 4 NAUTIKA / 31 -> DP 5 NAUTIKA / 38\n// Cell 3 (NO : 3): 3 -> T3, AMRUN SAIFUDDIN -> SAUFAN RUSDIONO,\n//                  DP 4 NAUTIKA / 31 -> DP 5 NAUTIKA / 38 (size stays 42)\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst row = table.rows.getFirst();\nconst cells = row.cells;\ncells.load(\"items\");\nawait context.sync();\n\n// Per-cell find/replace pairs, applied in document order of the cells.\nconst perCellReplacements = [\n  [\n    [\"1\", \"T1\"],\n    [\"ABU NAWAS\", \"GEDE SUDIANTARA\"],\n    [\"42\", \"43\"],\n    [\"DP 4 NAUTIKA / 31\", \"DP 5 NAUTIKA / 38\"],\n  ],\n  [\n    [\"2\", \"T2\"],\n    [\"ACHMAD UBAIDILLAH\", \"AGUS SETIAWAN\"],\n    [\"43\", \"42\"],\n    [\"DP 4 NAUTIKA / 31\", \"DP 5 NAUTIKA / 38\"],\n  ],\n  [\n    [\"3\", \"T3\"],\n    [\"AMRUN SAIFUDDIN\", \"SAUFAN RUSDIONO\"],\n    [\"DP 4 NAUTIKA / 31\", \"DP 5 NAUTIKA / 38\"],\n  ],\n];\n\n// `insertText(..., \"Replace\")` on a search hit whose new text ends with the\n// *same characters* as the matched text (e.g. \"1\" -> \"T1\") can land the\n// inserted prefix in the neighbouring run instead of the matched one. Route\n// every replacement through a throwaway placeholder that shares no\n// prefix/suffix with either the find or replace text to sidestep that.\nlet placeholderSeq = 0;\nasync function replaceInCellBody(cellBody, find, replace) {\n  const found = cellBody.search(find, { matchCase: true, matchWholeWord: true });\n  found.load(\"items\");\n  await context.sync();\n\n  for (let j = 0; j < found.items.length; j++) {\n    placeholderSeq++;\n    const placeholder = \"\\u0001PLACEHOLDER\" + placeholderSeq + \"\\u0001\";\n    found.items[j].insertText(placeholder, \"Replace\");\n    await context.sync();\n\n    const placeholderRange = cellBody.search(placeholder, { matchCase: true });\n    placeholderRange.load(\"items\");\n    await context.sync();\n    placeholderRange.items[0].insertText(replace, \"Replace\");\n    await context.sync();\n  }\n}\n\nfor (let i = 0; i < cells.items.length; i++) {\n  const cellBody = cells.items[i].body;\n  const replacements = perCellReplacements[i];\n  for (const [find, replace] of replacements) {\n    await replaceInCellBody(cellBody, find, replace);\n  }\n}\n", "ps1": "# Apply the per-label updates (NO / NAMA / SEPATU size / KELAS) to the\n# three shoe-label cells in the single-row table.\n#\n# Cell 1 (NO : 1): 1 -> T1, ABU NAWAS -> GEDE SUDIANTARA, 42 -> 43,\n#                  DP 4 NAUTIKA / 31 -> DP 5 NAUTIKA / 38\n# Cell 2 (NO : 2): 2 -> T2, ACHMAD UBAIDILLAH -> AGUS SETIAWAN, 43 -> 42,\n#                  DP 4 NAUTIKA / 31 -> DP 5 NAUTIKA / 38\n# Cell 3 (NO : 3): 3 -> T3, AMRUN SAIFUDDIN -> SAUFAN RUSDIONO,\n#                  DP 4 NAUTIKA / 31 -> DP 5 NAUTIKA / 38 (size stays 42)\n\nfunction Replace-InRange($range, $findText, $replaceText) {\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $true\n    # wdFindStop (0) so the search never leaves the supplied range, and\n    # wdReplaceOne (1) so only the single occurrence inside this cell's\n    # range is touched (wdReplaceAll here would leak across the whole\n    # document instead of staying scoped to $range).\n    $find.Execute($findText, $true, $true, $false, $false, $false, $true, 0, $false, $replaceText, 1) | Out-Null\n}\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n$cell1 = $tbl.Cell(1, 1)\n$cell2 = $tbl.Cell(1, 2)\n$cell3 = $tbl.Cell(1, 3)\n\nReplace-InRange $cell1.Range \"1\" \"T1\"\nReplace-InRange $cell1.Range \"ABU NAWAS\" \"GEDE SUDIANTARA\"\nReplace-InRange $cell1.Range \"42\" \"43\"\nReplace-InRange $cell1.Range \"DP 4 NAUTIKA / 31\" \"DP 5 NAUTIKA / 38\"\n\nReplace-InRange $cell2.Range \"2\" \"T2\"\nReplace-InRange $cell2.Range \"ACHMAD UBAIDILLAH\" \"AGUS SETIAWAN\"\nReplace-InRange $cell2.Range \"43\" \"42\"\nReplace-InRange $cell2.Range \"DP 4 NAUTIKA / 31\" \"DP 5 NAUTIKA / 38\"\n\nReplace-InRange $cell3.Range \"3\" \"T3\"\nReplace-InRange $cell3.Range \"AMRUN SAIFUDDIN\" \"SAUFAN RUSDIONO\"\nReplace-InRange $cell3.Range \"DP 4 NAUTIKA / 31\" \"DP 5 NAUTIKA / 38\"\n"}
